$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "Background`nid=`"Par1`">Frail older adults are predisposed to multiple comorbidities and adverse events.`n`n Recent interventional studies have shown that frailty can be improved and managed.`n`n In this study, effective individualized home-based exercise and nutrition interventions were developed for reducing frailty in older adults.`n`n`nMethods`nid=`"Par2`">This study was a four-arm, single-blind, randomized controlled trial conducted between October 2015 and June 2017 at Miaoli General Hospital in Taiwan.`n`n Overall, 319 pre-frail or frail older adults were randomly assigned into one of the four study groups (control, exercise, nutrition, and exercise plus nutrition [combination]) and followed up during a 3-month intervention period and 3-month self-maintenance period.`n`n Improvement in frailty scores was the primary outcome.`n`n Secondary outcomes included improvements in physical performance and mental health.`n`n The measurements were performed at baseline, 1 month, 3 months, and 6 months.`n`n`nResults`nid=`"Par3`">At the 6-month measurement, the exercise (difference in frailty score change from baseline: − 0.23; 95% confidence interval [CI]: − 0.41, − 0.05; p = 0.012), nutrition (− 0.28; 95% CI: − 0.46, − 0.11; p = 0.002), and combination (− 0.34; 95% CI: − 0.52, − 0.16; p &lt;  0.001) groups exhibited significantly greater improvements in the frailty scores than the control group.`n`n Significant improvements were also observed in several physical performance parameters in the exercise, nutrition, and combination groups, as well as in the 12-Item Short Form Health Survey mental component summary score for the nutrition group.`n`n`nConclusions`nid=`"Par4`">The designated home-based exercise and nutrition interventions can help pre-frail or frail older adults to improve their frailty score and physical performance.`n`n`nTrial registration`nid=`"Par5`">Retrospectively registered at ClinicalTrials.`n`ngov (identifier: NCT03477097); registration date: March 26, 2018.`n"
$ws.Range("E5").Value = "[Tsung-Jen%Hsieh%tjhsieh@nhri.edu.tw%1,    Shin-Chang%Su%maxmax0220@hotmail.com%1,    Chun-Wei%Chen%chunwei.peter@gmail.com%1,    Yaw-Wen%Kang%castle.kang@gmail.com%1,    Ming-Hsia%Hu%mhh@ntu.edu.tw%1,    Li-Lin%Hsu%llh@firdi.org.tw%1,    Szu-Yun%Wu%s.wu@ibms.sinica.edu.tw%1,    Likwang%Chen%likwang@nhri.edu.tw%1,    Hsing-Yi%Chang%hsingyi@nhri.edu.tw%1,    Shao-Yuan%Chuang%chuangsy@nhri.edu.tw%1,    Wen-Harn%Pan%pan@ibms.sinica.edu.tw%1,    Chih-Cheng%Hsu%cch@nhri.edu.tw%1]"
